$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells need Text number format applied first so that
# numeric-looking strings (e.g. "242.82", "1.000", "0.000007726") keep
# their exact original text representation instead of being coerced to
# floating point numbers by Excel's normal auto-detection. NumberFormat is
# applied per contiguous block (union/comma ranges only honor the first area).
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31:D35").NumberFormat = "@"
$ws.Range("D38:D39").NumberFormat = "@"
$ws.Range("D41:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '30.810.64'
$ws.Cells.Item(2, 5).Value = '  -1.04%  '
$ws.Cells.Item(3, 4).Value = '1.941.97'
$ws.Cells.Item(3, 5).Value = '  -0.74%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).Value = '242.82'
$ws.Cells.Item(5, 5).Value = '  -1.56%  '
$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 5).Value = '  -0.03%  '
$ws.Cells.Item(7, 4).Value = '0.4883'
$ws.Cells.Item(7, 5).Value = '  -0.08%  '
$ws.Cells.Item(8, 4).Value = '0.2947'
$ws.Cells.Item(8, 5).Value = '  -0.59%  '
$ws.Cells.Item(9, 4).Value = '0.06891'
$ws.Cells.Item(9, 5).Value = '  +0.85%  '
$ws.Cells.Item(10, 4).Value = '19.44'
$ws.Cells.Item(10, 5).Value = '  +2.00%  '
$ws.Cells.Item(11, 4).Value = '105.97'
$ws.Cells.Item(11, 5).Value = '  -0.42%  '
$ws.Cells.Item(12, 4).Value = '1.941.02'
$ws.Cells.Item(12, 5).Value = '  -0.06%  '
$ws.Cells.Item(13, 4).Value = '0.07732'
$ws.Cells.Item(13, 5).Value = '  -0.16%  '
$ws.Cells.Item(14, 4).Value = '5.358'
$ws.Cells.Item(14, 5).Value = '  -0.96%  '
$ws.Cells.Item(15, 4).Value = '0.6988'
$ws.Cells.Item(15, 5).Value = '  -2.25%  '
$ws.Cells.Item(16, 4).Value = '273.16'
$ws.Cells.Item(16, 5).Value = '  -3.89%  '
$ws.Cells.Item(17, 4).Value = '30.810.37'
$ws.Cells.Item(17, 5).Value = '  -0.78%  '
$ws.Cells.Item(18, 4).Value = '0.000007726'
$ws.Cells.Item(18, 5).Value = '  -0.25%  '
$ws.Cells.Item(19, 4).Value = '13.11'
$ws.Cells.Item(19, 5).Value = '  -0.78%  '
$ws.Cells.Item(20, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(20, 4).Value = '2.197.54'
$ws.Cells.Item(20, 5).Value = '  +0.08%  '
$ws.Cells.Item(21, 2).Value = 'Dai'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(21, 4).Value = '1.001'
$ws.Cells.Item(21, 5).Value = '  -0.01%  '
$ws.Cells.Item(22, 4).Value = '5.508'
$ws.Cells.Item(22, 5).Value = '  -0.29%  '
$ws.Cells.Item(23, 5).Value = '  -0.18%  '
$ws.Cells.Item(24, 4).Value = '6.565'
$ws.Cells.Item(24, 5).Value = '  -0.57%  '
$ws.Cells.Item(26, 4).Value = '167.13'
$ws.Cells.Item(26, 5).Value = '  -1.07%  '
$ws.Cells.Item(27, 5).Value = '  -1.27%  '
$ws.Cells.Item(28, 5).Value = '  -1.83%  '
$ws.Cells.Item(29, 4).Value = '0.1043'
$ws.Cells.Item(29, 5).Value = '  -0.70%  '
$ws.Cells.Item(31, 4).Value = '4.576'
$ws.Cells.Item(31, 5).Value = '  -3.32%  '
$ws.Cells.Item(32, 4).Value = '1.554'
$ws.Cells.Item(33, 4).Value = '4.376'
$ws.Cells.Item(33, 5).Value = '  -2.71%  '
$ws.Cells.Item(34, 4).Value = '0.04861'
$ws.Cells.Item(34, 5).Value = '  -2.62%  '
$ws.Cells.Item(35, 4).Value = '0.7548'
$ws.Cells.Item(35, 5).Value = '  -0.99%  '
$ws.Cells.Item(36, 5).Value = '  -0.75%  '
$ws.Cells.Item(37, 5).Value = '  +0.06%  '
$ws.Cells.Item(38, 4).Value = '2.734'
$ws.Cells.Item(38, 5).Value = '  +0.05%  '
$ws.Cells.Item(39, 4).Value = '0.01998'
$ws.Cells.Item(39, 5).Value = '  -2.45%  '
$ws.Cells.Item(40, 5).Value = '  -2.20%  '
$ws.Cells.Item(41, 4).Value = '6.557'
$ws.Cells.Item(41, 5).Value = '  +1.80%  '
$ws.Cells.Item(42, 4).Value = '77.91'
$ws.Cells.Item(42, 5).Value = '  +7.21%  '
$ws.Cells.Item(43, 4).Value = '2.100'
$ws.Cells.Item(43, 5).Value = '  -2.38%  '
$ws.Cells.Item(44, 4).Value = '0.9038'
$ws.Cells.Item(44, 5).Value = '  +2.28%  '
$ws.Cells.Item(45, 4).Value = '108.14'
$ws.Cells.Item(45, 5).Value = '  -1.47%  '
$ws.Cells.Item(46, 4).Value = '0.4406'
$ws.Cells.Item(46, 5).Value = '  -1.18%  '
$ws.Cells.Item(47, 4).Value = '0.9992'
$ws.Cells.Item(47, 5).Value = '  -0.08%  '
$ws.Cells.Item(48, 4).Value = '7.762'
$ws.Cells.Item(48, 5).Value = '  +3.33%  '
$ws.Cells.Item(49, 4).Value = '1.005.35'
$ws.Cells.Item(49, 5).Value = '  +3.23%  '
$ws.Cells.Item(50, 4).Value = '0.1247'
$ws.Cells.Item(51, 4).Value = '9.288'
$ws.Cells.Item(51, 5).Value = '  -0.94%  '
